$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for columns B,D,E,F,G,H,I,J,K,L, rows 2-25
# Column order per row: B, D, E, F, G, H, I, J, K, L
$data = @{
    2 = @(14.03311277436778, 7.388192954733805, 12.87265837027168, 40.58806043070942, 47.91711018936709, 19.06477042075052, 30.70501778969327, 10.40892301417651, 11.8839660152618, 11.0103834114437)
    3 = @(13.92420072254983, 7.38961780345785, 12.88998932666533, 40.6408985370179, 47.9270237027513, 19.10932683946282, 30.80580600551322, 10.42161312577418, 11.67134503735698, 10.98397211551211)
    4 = @(13.85942672162227, 7.390960084707104, 12.90136481853177, 40.68203147786414, 47.94583086333256, 19.13988751981443, 30.8724136939007, 10.42983606565645, 11.54129312456767, 10.96954968075215)
    5 = @(13.83358292304378, 7.391625054949845, 12.9061854763464, 40.70097518542708, 47.95668890611893, 19.15314571474011, 30.90074415273279, 10.43329569576371, 11.48849356959087, 10.96412775721098)
    6 = @(13.82932562429163, 7.391742611278216, 12.90699713224924, 40.70425242704609, 47.95868460942661, 19.15539578786076, 30.90552010059873, 10.43387673958373, 11.47974031033907, 10.96325506667277)
    7 = @(13.85907591591728, 7.390968574436934, 12.90142908175498, 40.68227813090463, 47.94596437410794, 19.1400630684881, 30.87279096169562, 10.42988228285596, 11.54058015486486, 10.96947470993035)
    8 = @(13.99513931337512, 7.388587497083436, 12.87848195510916, 40.60447346348327, 47.91788608770037, 19.07946818844306, 30.7387889841901, 10.41320925621795, 11.81059404964571, 11.00090701551052)
    9 = @(14.27749500523987, 7.38760776333276, 12.83928923201478, 40.52098201116154, 47.96387642840126, 18.98609432988077, 30.51350795702251, 10.38392102719019, 12.34091782841432, 11.07658393988806)
    10 = @(14.49287277089721, 7.38911181675891, 12.81400793544475, 40.50189662232489, 48.05931275086711, 18.9330602479775, 30.37087857132837, 10.36446133830929, 12.7269817407617, 11.14045547301791)
    11 = @(14.5922377579708, 7.390273412838609, 12.80326424169993, 40.50240421153173, 48.11607998658808, 18.91232361747078, 30.31096785090297, 10.35605155742546, 12.90102150682047, 11.17124029109452)
    12 = @(14.63003552472333, 7.390781433370631, 12.79930429997163, 40.50391758398438, 48.1394895449462, 18.90495915487586, 30.28899682300623, 10.35293032495925, 12.96663401622765, 11.18314007086354)
    13 = @(14.6218879874965, 7.390668999160088, 12.80015232639549, 40.50353290696811, 48.13436293141723, 18.90652350672091, 30.29369682470167, 10.35359972407315, 12.95251720062132, 11.18056657260074)
    14 = @(14.59534415580118, 7.390313845074107, 12.80293628332712, 40.50250224430846, 48.11796761344443, 18.91170795355366, 30.30914592920369, 10.35579350325713, 12.90642571169991, 11.17221447631803)
    15 = @(14.57910661376231, 7.390105163006891, 12.80465565110508, 40.50204296434388, 48.10817389338033, 18.91494715576133, 30.31870220441274, 10.35714550101746, 12.87815336902664, 11.16712992625762)
    16 = @(14.4864044769434, 7.389045467649994, 12.81472525930617, 40.50204836279809, 48.05587101722727, 18.93448369279242, 30.37489400412174, 10.36501981927889, 12.71557035664135, 11.1384778680999)
    17 = @(14.42986957481514, 7.388517298767634, 12.82109622610513, 40.50440541450806, 48.0272005457235, 18.94733724183428, 30.4106399498993, 10.36996361064898, 12.61537833455364, 11.12133920118813)
    18 = @(14.39748465529278, 7.388258475167525, 12.82483190267231, 40.50662617048758, 48.01196744174311, 18.95504919811971, 30.43166801681065, 10.37284881738403, 12.55760525797407, 11.11164467235387)
    19 = @(14.38654331829708, 7.388178579801543, 12.82610898998383, 40.50752664746145, 48.00702591160304, 18.95771508747029, 30.43886810228896, 10.37383286359741, 12.53802136865034, 11.10839048871924)
    20 = @(14.43587431416853, 7.388568873120985, 12.82041065292256, 40.50406497313793, 48.03012247492494, 18.94593594537919, 30.40678629609855, 10.36943302500249, 12.62605946285479, 11.12314679877289)
    21 = @(14.60313633157905, 7.390416316915226, 12.80211562706674, 40.50276912559467, 48.122731464003, 18.91017190488214, 30.30458871737128, 10.35514741942352, 12.91997234135819, 11.17466116804514)
    22 = @(14.7134313636732, 7.392020716006607, 12.79079079610991, 40.50962220479573, 48.19440257277495, 18.88964304116811, 30.24196969775351, 10.34618019194909, 13.11032825813641, 11.2097375308517)
    23 = @(14.65448487607193, 7.391128261090202, 12.79677736521384, 40.5052603554468, 48.15513343242776, 18.90033914322374, 30.27500852321687, 10.350932471636, 13.00891092327281, 11.19088992145202)
    24 = @(14.4331592025009, 7.388545416709783, 12.82072037340419, 40.5042161901377, 48.02879757667537, 18.94656846763549, 30.4085270464024, 10.36967276895639, 12.62123105346823, 11.12232908901962)
    25 = @(14.19960382641056, 7.387480400522503, 12.8492729975108, 40.53615523331097, 47.94060348844324, 19.00862398474359, 30.57043559026915, 10.39148143805495, 12.19779102909269, 11.05463659283838)
}

$cols = @("B","D","E","F","G","H","I","J","K","L")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}
